$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold text-formatted values (e.g. "301.85", "-0.96%").
# Force text format on just the cells being updated so Excel does not
# auto-convert the assigned strings into numbers/percentages.
$changedCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","E25","E26","D27","E27","D39","E39","D40","E40","D41","E41","E42","D43","E43","D44","E44","D45","E45","D46","E46","E47","D49","E49","E50","E51")
foreach ($addr in $changedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "301.85"
$ws.Range("E2").Value = "-0.96%"
$ws.Range("D3").Value = "37.41"
$ws.Range("E3").Value = "5.67%"
$ws.Range("D4").Value = "5.003"
$ws.Range("E4").Value = "-2.82%"
$ws.Range("D5").Value = "0.07844"
$ws.Range("E5").Value = "0.77%"
$ws.Range("D6").Value = "2.228"
$ws.Range("E6").Value = "-7.74%"
$ws.Range("D7").Value = "8.031"
$ws.Range("E7").Value = "0.05%"
$ws.Range("E8").Value = "1.99%"
$ws.Range("D9").Value = "0.9090"
$ws.Range("E9").Value = "-1.50%"
$ws.Range("D10").Value = "0.09564"
$ws.Range("E10").Value = "-2.96%"
$ws.Range("D11").Value = "0.1890"
$ws.Range("E11").Value = "4.98%"
$ws.Range("D12").Value = "0.08482"
$ws.Range("E12").Value = "-1.35%"
$ws.Range("D13").Value = "0.03525"
$ws.Range("E13").Value = "6.26%"
$ws.Range("D14").Value = "0.09955"
$ws.Range("E14").Value = "0.41%"
$ws.Range("D15").Value = "0.001484"
$ws.Range("E15").Value = "-0.09%"
$ws.Range("D16").Value = "0.005694"
$ws.Range("E16").Value = "-0.22%"
$ws.Range("D17").Value = "3.465"
$ws.Range("E17").Value = "-0.17%"
$ws.Range("D18").Value = "2.069"
$ws.Range("E18").Value = "-3.35%"
$ws.Range("E19").Value = "2.85%"
$ws.Range("D20").Value = "0.1309"
$ws.Range("E20").Value = "1.22%"
$ws.Range("D21").Value = "4.760"
$ws.Range("E21").Value = "10.13%"
$ws.Range("D22").Value = "0.2204"
$ws.Range("E22").Value = "2.68%"
$ws.Range("D23").Value = "0.04650"
$ws.Range("E23").Value = "1.81%"
$ws.Range("D24").Value = "0.001229"
$ws.Range("E24").Value = "0.98%"
$ws.Range("E25").Value = "-0.02%"
$ws.Range("E26").Value = "0.07%"
$ws.Range("D27").Value = "0.0004747"
$ws.Range("E27").Value = "28.36%"
$ws.Range("D39").Value = "0.01759"
$ws.Range("E39").Value = "-1.88%"
$ws.Range("D40").Value = "0.04744"
$ws.Range("E40").Value = "-0.06%"
$ws.Range("D41").Value = "0.007800"
$ws.Range("E41").Value = "0.60%"
$ws.Range("E42").Value = "-1.51%"
$ws.Range("D43").Value = "0.007662"
$ws.Range("E43").Value = "7.26%"
$ws.Range("D44").Value = "0.002229"
$ws.Range("E44").Value = "3.66%"
$ws.Range("D45").Value = "0.009820"
$ws.Range("E45").Value = "2.82%"
$ws.Range("D46").Value = "0.00006073"
$ws.Range("E46").Value = "-0.72%"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("D49").Value = "0.002688"
$ws.Range("E49").Value = "34.47%"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("E51").Value = "-0.01%"
